$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 21000000.0
$ws.Range("C12").Value = 21000000.0
$ws.Range("D12").Value = 18000000.0
$ws.Range("E12").Value = 14000000.0
$ws.Range("F12").Value = 15000000.0
